$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 284.77777
$ws.Range("I9").Value = 211.85715
$ws.Range("K9").Value = 211.85715
$ws.Range("M9").Value = -42.85714999999999

$ws.Range("H11").Value = 42.636364
$ws.Range("I11").Value = 42.636364
$ws.Range("K11").Value = 42.636364
$ws.Range("M11").Value = 97.363636

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H40").Value = 9028.857
$ws.Range("J40").Value = 11200.4
$ws.Range("L40").Value = 11200.4
$ws.Range("N40").Value = -11550.4

$ws.Range("H55").Value = 677.61536
$ws.Range("I55").Value = 405.45456
$ws.Range("K55").Value = 405.45456
$ws.Range("M55").Value = -191.45456

$ws.Range("H113").Value = 11500.25
$ws.Range("J113").Value = 19001
$ws.Range("L113").Value = 19001
$ws.Range("N113").Value = -25509

$ws.Range("H118").Value = 502
$ws.Range("I118").Value = 502
$ws.Range("K118").Value = 1506
$ws.Range("M118").Value = 151

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 1257.697
$ws.Range("I132").Value = 1079.7037
$ws.Range("K132").Value = 3239.1111
$ws.Range("M132").Value = -709.1111000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10355.454
$ws.Range("I2").Value = 1022.8571
$ws.Range("J2").Value = 26687.5
$ws.Range("K2").Value = 1022.8571
$ws.Range("L2").Value = 26687.5
$ws.Range("M2").Value = -909.8570999999999
$ws.Range("N2").Value = -26913.5

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H45").Value = 45457390
$ws.Range("I45").Value = 71429870
$ws.Range("K45").Value = 71429870
$ws.Range("M45").Value = -71429493

$ws.Range("H74").Value = 47625164
$ws.Range("I74").Value = 55561508
$ws.Range("K74").Value = 55561508
$ws.Range("M74").Value = -55560634

$ws.Range("H77").Value = 47625164
$ws.Range("I77").Value = 55561508
$ws.Range("K77").Value = 277807540
$ws.Range("M77").Value = -277803172

$ws.Range("H88").Value = 2450
$ws.Range("J88").Value = 2412.7144
$ws.Range("L88").Value = 2412.7144
$ws.Range("N88").Value = -3224.7144

$ws.Range("H91").Value = 2450
$ws.Range("J91").Value = 2412.7144
$ws.Range("L91").Value = 2412.7144
$ws.Range("N91").Value = -5220.7144

$ws.Range("H101").Value = 25301
$ws.Range("J101").Value = 25301
$ws.Range("L101").Value = 25301
$ws.Range("N101").Value = -31791

$ws.Range("H116").Value = 10355.454
$ws.Range("I116").Value = 1022.8571
$ws.Range("J116").Value = 26687.5
$ws.Range("K116").Value = 1022.8571
$ws.Range("L116").Value = 26687.5
$ws.Range("M116").Value = 1271.1429
$ws.Range("N116").Value = -31275.5

$ws.Range("H130").Value = 44476
$ws.Range("J130").Value = 44476
$ws.Range("L130").Value = 44476
$ws.Range("N130").Value = -54516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10355.454
$ws.Range("I3").Value = 1022.8571
$ws.Range("J3").Value = 26687.5
$ws.Range("K3").Value = 1022.8571
$ws.Range("L3").Value = 26687.5
$ws.Range("M3").Value = -908.8570999999999
$ws.Range("N3").Value = -26915.5

$ws.Range("H86").Value = 2214.2285
$ws.Range("I86").Value = 1614.258
$ws.Range("K86").Value = 1614.258
$ws.Range("M86").Value = -491.258

$ws.Range("H89").Value = 2214.2285
$ws.Range("I89").Value = 1614.258
$ws.Range("K89").Value = 8071.29
$ws.Range("M89").Value = -2455.29

$ws.Range("H95").Value = 33998
$ws.Range("J95").Value = 33998
$ws.Range("L95").Value = 33998
$ws.Range("N95").Value = -39490

$ws.Range("H134").Value = 5405.7
$ws.Range("I134").Value = 1823.2
$ws.Range("J134").Value = 8988.200000000001
$ws.Range("K134").Value = 5469.6
$ws.Range("L134").Value = 26964.6
$ws.Range("M134").Value = -2934.6
$ws.Range("N134").Value = -32034.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 27335.666
$ws.Range("I17").Value = 40504
$ws.Range("J17").Value = 999
$ws.Range("K17").Value = 40504
$ws.Range("L17").Value = 999
$ws.Range("M17").Value = -40330
$ws.Range("N17").Value = -1347

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H28").Value = 48475
$ws.Range("J28").Value = 48475
$ws.Range("L28").Value = 48475
$ws.Range("N28").Value = -48965

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 8338.5
$ws.Range("I122").Value = 4102.125
$ws.Range("J122").Value = 16811.25
$ws.Range("K122").Value = 12306.375
$ws.Range("L122").Value = 50433.75
$ws.Range("M122").Value = -9856.375
$ws.Range("N122").Value = -55333.75

$ws.Range("H132").Value = 6032.75
$ws.Range("J132").Value = 7064.25
$ws.Range("L132").Value = 21192.75
$ws.Range("N132").Value = -26252.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2119.4
$ws.Range("I113").Value = 1866
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 5598
$ws.Range("L113").Value = 7498.5
$ws.Range("M113").Value = -3428
$ws.Range("N113").Value = -11838.5

$ws.Range("H132").Value = 5252.8667
$ws.Range("I132").Value = 5073.75
$ws.Range("J132").Value = 5318
$ws.Range("K132").Value = 45663.75
$ws.Range("L132").Value = 47862
$ws.Range("M132").Value = -43133.75
$ws.Range("N132").Value = -52922

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9408.223
$ws.Range("I22").Value = 1124.6666
$ws.Range("K22").Value = 1124.6666
$ws.Range("M22").Value = -829.6666

$ws.Range("H27").Value = 9408.223
$ws.Range("I27").Value = 1124.6666
$ws.Range("K27").Value = 1124.6666
$ws.Range("M27").Value = -1017.6666

$ws.Range("H100").Value = 3659.9644
$ws.Range("I100").Value = 2590.0454
$ws.Range("K100").Value = 2590.0454
$ws.Range("M100").Value = -2049.0454

$ws.Range("H128").Value = 99975
$ws.Range("J128").Value = 99975
$ws.Range("L128").Value = 99975
$ws.Range("N128").Value = -109935

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 35000
$ws.Range("J82").Value = 35000
$ws.Range("L82").Value = 35000
$ws.Range("N82").Value = -35766

$ws.Range("H85").Value = 35000
$ws.Range("J85").Value = 35000
$ws.Range("L85").Value = 35000
$ws.Range("N85").Value = -37652

$ws.Range("H112").Value = 40155
$ws.Range("J112").Value = 40155
$ws.Range("L112").Value = 40155
$ws.Range("N112").Value = -43109

$ws.Range("H124").Value = 58806.668
$ws.Range("J124").Value = 58806.668
$ws.Range("L124").Value = 58806.668
$ws.Range("N124").Value = -68626.66800000001

$ws.Range("H132").Value = 5749.5
$ws.Range("I132").Value = 5771.273
$ws.Range("J132").Value = 5701.6
$ws.Range("K132").Value = 17313.819
$ws.Range("L132").Value = 17104.8
$ws.Range("M132").Value = -14783.819
$ws.Range("N132").Value = -22164.8

$ws.Range("H136").Value = 5527.769
$ws.Range("I136").Value = 3395.85
$ws.Range("K136").Value = 10187.55
$ws.Range("M136").Value = -7637.549999999999

Write-Host "edit.ps1 applied: 43 leve rows updated across 8 sheets"
